$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.867.11"
$ws.Range("E2").Value = "  -4.76%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.219.93"
$ws.Range("E3").Value = "  -6.13%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.92"
$ws.Range("E5").Value = "  +1.78%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "99.67"
$ws.Range("E6").Value = "  -7.64%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.591"
$ws.Range("E7").Value = "  -6.11%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.563"
$ws.Range("E9").Value = "  -7.83%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.12"
$ws.Range("E10").Value = "  -8.72%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.88"
$ws.Range("E11").Value = "  -3.01%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0828"
$ws.Range("E12").Value = "  -9.47%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.81"
$ws.Range("E13").Value = "  -7.14%  "
$ws.Range("E14").Value = "  -3.25%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.863"
$ws.Range("E15").Value = "  -11.23%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.562.14"
$ws.Range("E16").Value = "  -5.99%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.25"
$ws.Range("E17").Value = "  -6.09%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.210.66"
$ws.Range("E18").Value = "  -6.86%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "42.872.15"
$ws.Range("E19").Value = "  -4.80%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.48"
$ws.Range("E20").Value = "  +7.97%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0963"
$ws.Range("E21").Value = "  -8.81%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.44"
$ws.Range("E22").Value = "  -10.83%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.37"
$ws.Range("E24").Value = "  -8.76%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "236.45"
$ws.Range("E25").Value = "  -8.65%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.13"
$ws.Range("E26").Value = "  -7.20%  "
$ws.Range("E27").Value = "  -0.19%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.11"
$ws.Range("E28").Value = "  -9.03%  "
$ws.Range("E29").Value = "  -5.18%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.35"
$ws.Range("E30").Value = "  -11.27%  "
$ws.Range("E31").Value = "  -7.96%  "
$ws.Range("E32").Value = "  -7.80%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "34.31"
$ws.Range("E33").Value = "  -7.67%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "157.23"
$ws.Range("E34").Value = "  -6.68%  "
$ws.Range("E35").Value = "  -6.19%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.22"
$ws.Range("E36").Value = "  +10.45%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.97"
$ws.Range("E37").Value = "  +13.30%  "
$ws.Range("E38").Value = "  -5.45%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.44"
$ws.Range("E39").Value = "  -5.69%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.79"
$ws.Range("E40").Value = "  -2.83%  "
$ws.Range("E41").Value = "  -11.07%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0325"
$ws.Range("E42").Value = "  -7.41%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.915.59"
$ws.Range("E43").Value = "  +1.14%  "
$ws.Range("E44").Value = "  -0.15%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.33"
$ws.Range("E45").Value = "  -3.51%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "89.20"
$ws.Range("E46").Value = "  -10.45%  "
$ws.Range("E47").Value = "  -8.97%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.41"
$ws.Range("E48").Value = "  -3.47%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "60.70"
$ws.Range("E49").Value = "  -12.49%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "74.85"
$ws.Range("E50").Value = "  -7.20%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.865"
$ws.Range("E51").Value = "  +17.67%  "
